# Weekly data refresh: insert a new observation at row 11 (pushing the
# existing rows 11:128 down to 12:129) and populate the new row with this
# week's figures for "Comercializadora del Agro de Limarí" / Poroto granado.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data rows 11:128 down to 12:129, growing the used range to
# A1:R129 and creating a blank row 11 ready to be filled in.
$ws.Rows("11:11").Insert()

# Populate the newly inserted row 11 with the new week's record.
$ws.Range("A11").Value2 = 2
$ws.Range("B11").Value2 = 'Comercializadora del Agro de Limarí'
$ws.Range("C11").Value2 = 'Coquimbo'
$ws.Range("D11").Value2 = 45050
$ws.Range("E11").Value2 = 4
$ws.Range("F11").Value2 = 100112030
$ws.Range("G11").Value2 = 'Poroto granado'
$ws.Range("H11").Value2 = 'Sin especificar'
$ws.Range("I11").Value2 = 'Primera'
$ws.Range("J11").Value2 = 600
$ws.Range("K11").Value2 = 29000
$ws.Range("L11").Value2 = 30000
$ws.Range("M11").Value2 = 29500
$ws.Range("N11").Value2 = '$/malla 25 kilos'
$ws.Range("O11").Value2 = 'Provincia de Limarí'
$ws.Range("P11").Value2 = 1180
$ws.Range("Q11").Value2 = 25
$ws.Range("R11").Value2 = 'Hortaliza'
